# Add the new CQ4-3 competency question (for the "recipe video" CQ) to the
# "Competency Questions" sheet, right after the existing CQ4-2 row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Competency Questions")
$ws.Activate()

# Row 18 is the last CQ4 row (CQ4-2). Insert a new row below it so the
# REQ4 group gets a third competency question, and everything that used to
# follow (REQ5..REQ7 blocks) shifts down by one row.
$ws.Rows.Item(19).Insert()
$ws.Rows.Item(19).RowHeight = 43.15

$ws.Cells.Item(19, 1).Value = "REQ4"
$ws.Cells.Item(19, 3).Value = "Where is it possible to find a video showing how to prepare recipe Z?"
$ws.Cells.Item(19, 2).Value = "CQ4-3"

# Restore the view state (scroll position / selection) recorded after the edit.
$ws.Range("F21").Select()
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1
